# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.298.29"
$ws.Range("E2").Value = "  -4.67%  "
$ws.Range("D3").Value = "2.953.89"
$ws.Range("E3").Value = "  -6.53%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'536.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.59%  "
$ws.Range("D6").Value = "'151.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.38%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.560"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.32%  "
$ws.Range("D9").Value = "2.959.16"
$ws.Range("E9").Value = "  -6.21%  "
$ws.Range("E10").Value = "  -4.79%  "
$ws.Range("E11").Value = "  -8.24%  "
$ws.Range("D12").Value = "'0.362"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.69%  "
$ws.Range("D13").Value = "3.469.46"
$ws.Range("E13").Value = "  -6.44%  "
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "61.358.31"
$ws.Range("E15").Value = "  -4.67%  "
$ws.Range("D16").Value = "'23.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.36%  "
$ws.Range("D17").Value = "2.957.67"
$ws.Range("E17").Value = "  -6.23%  "
$ws.Range("E18").Value = "  -6.57%  "
$ws.Range("D19").Value = "'5.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("D20").Value = "'380.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.24%  "
$ws.Range("E21").Value = "  -6.36%  "
$ws.Range("E22").Value = "  -6.57%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'64.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.51%  "
$ws.Range("D25").Value = "3.078.81"
$ws.Range("E25").Value = "  -6.81%  "
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("D27").Value = "'0.185"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.00%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  -10.65%  "
$ws.Range("D30").Value = "'8.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.99%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  -6.23%  "
$ws.Range("D33").Value = "'20.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.56%  "
$ws.Range("D34").Value = "'158.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").Value = "'5.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.22%  "
$ws.Range("D36").Value = "'4.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.88%  "
$ws.Range("E37").Value = "  -6.97%  "
$ws.Range("E38").Value = "  -6.09%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.13%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.401.54"
$ws.Range("E40").Value = "  -10.75%  "
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'36.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.68%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'22.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.71%  "
$ws.Range("D44").Value = "'0.660"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("E45").Value = "  -5.34%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -5.88%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.62%  "
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("D51").Value = "'19.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.88%  "
